# "Repaired sheets naming error"
#
# The commit fixes the worksheet tab name (which had lost its space,
# "Sheet1") back to "Sheet 1", and also relabels the built-in "Normal"
# cell style using its localized Polish name "Normalny".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Sheet1" -> "Sheet 1"
$ws.Name = "Sheet 1"

# Rename the built-in cell style: "Normal" -> "Normalny"
$wb.Styles.Item(1).Name = "Normalny"
